# Add formatting to pivot example:
#  - bold the "Stock Value Summary" / "Category" header row of the pivot block
#  - insert a "Clothing" row into the pivot summary (rows shift up: Clothing/Stickers/Pets)
#  - best-fit the column widths for the whole used range (A:E)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the pivot summary rows (12-16 -> 12-15) -----------------------
# Row 16 (old "Pets" row) becomes row 15; clear its old location first.
$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()

$ws.Range("A15").Value = "Pets"
$ws.Range("B15").Formula = '=SUMIF(B2:B7, "Pets", E2:E7)'

$ws.Range("A14").Value = "Stickers"
$ws.Range("B14").Formula = '=SUMIF(B2:B7, "Stickers", E2:E7)'

# Newly inserted row for "Clothing"
$ws.Range("A13").Value = "Clothing"
$ws.Range("B13").Formula = '=SUMIF(B2:B7, "Clothing", E2:E7)'

# --- Bold the pivot header row ----------------------------------------------
$ws.Range("A12:B12").Font.Bold = $true

# --- Best-fit the column widths ---------------------------------------------
$ws.Columns("A:E").AutoFit()
$ws.Columns("A").ColumnWidth = 30.583333333333332
$ws.Columns("B").ColumnWidth = 10.25
$ws.Columns("C").ColumnWidth = 6.916666666666667
$ws.Columns("D").ColumnWidth = 10.416666666666666
$ws.Columns("E").ColumnWidth = 6.75

Write-Host "done"
